$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the duplicated rows (3-6) and the duplicate summary row (8) entirely,
# without shifting remaining rows up.
$ws.Rows.Item(3).Clear()
$ws.Rows.Item(4).Clear()
$ws.Rows.Item(5).Clear()
$ws.Rows.Item(6).Clear()
$ws.Rows.Item(8).Clear()

# Update the selected cell to match the new view state
$ws.Range("S17").Select()
